$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 26
$ws.Range("F4").Value = 606
$ws.Range("F5").Value = 32
$ws.Range("F6").Value = 1946
$ws.Range("F7").Value = 5479
$ws.Range("F8").Value = 1532
$ws.Range("F10").Value = 3136
$ws.Range("F11").Value = 696
$ws.Range("F13").Value = 1297
$ws.Range("F14").Value = 4344
$ws.Range("F15").Value = 1037
$ws.Range("F16").Value = 1679
$ws.Range("F17").Value = 2608
$ws.Range("F18").Value = 35
$ws.Range("F19").Value = 35
$ws.Range("F20").Value = 141
$ws.Range("F21").Value = 151
$ws.Range("F22").Value = 991
$ws.Range("F23").Value = 297
$ws.Range("F24").Value = 77
$ws.Range("F25").Value = 2
$ws.Range("F28").Value = 1097
$ws.Range("F29").Value = 392
$ws.Range("F30").Value = 60
$ws.Range("F31").Value = 177
$ws.Range("F32").Value = 309
$ws.Range("F34").Value = 4
$ws.Range("F35").Value = 1690
$ws.Range("F36").Value = 2203
$ws.Range("F37").Value = 1027
$ws.Range("F39").Value = 254
$ws.Range("F40").Value = 614
$ws.Range("F41").Value = 307
$ws.Range("F42").Value = 3
$ws.Range("F43").Value = 6
$ws.Range("F46").Value = 408
$ws.Range("F47").Value = 343

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 35
$ws.Range("F10").Value = 150
$ws.Range("F11").Value = 5

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 751

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 751
$ws.Range("F4").Value = 26
$ws.Range("F6").Value = 606
$ws.Range("F7").Value = 32
$ws.Range("F8").Value = 1946
$ws.Range("F9").Value = 5479
$ws.Range("F10").Value = 1532
$ws.Range("F13").Value = 3136
$ws.Range("F15").Value = 1297
$ws.Range("F16").Value = 4344
$ws.Range("F17").Value = 1037
$ws.Range("F18").Value = 1679
$ws.Range("F19").Value = 35
$ws.Range("F20").Value = 35
$ws.Range("F23").Value = 35
$ws.Range("F24").Value = 151
$ws.Range("F25").Value = 150
$ws.Range("F26").Value = 991
$ws.Range("F27").Value = 297
$ws.Range("F30").Value = 1097
$ws.Range("F31").Value = 392
$ws.Range("F32").Value = 60
$ws.Range("F33").Value = 177
$ws.Range("F35").Value = 1690
$ws.Range("F36").Value = 2203
$ws.Range("F37").Value = 1027
$ws.Range("F41").Value = 254
$ws.Range("F42").Value = 614
$ws.Range("F43").Value = 307
$ws.Range("F45").Value = 408
$ws.Range("F46").Value = 343

